# Applies the TestResultsInOneFile.xlsx edit described by the commit:
# "Edited the test files, for some of them were not valid"
#
# Semantic changes on Sheet1:
#  - Row for "CPSC313": test now succeeds (was a hard ERROR/NO VALID SOLUTION);
#    column B gets the Eval-value result, column C (console-error) is cleared.
#  - Row for "CPSC413": same kind of fix.
#  - Row for "EveningClass": the recorded result (column E) is corrected to a
#    result that actually satisfies the constraint, and the memo explaining
#    why the old result was wrong (column F) is removed; the row's yellow
#    "needs attention" highlight is cleared since it is no longer a problem
#    case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: CPSC313 -------------------------------------------------
$ws.Range("B14").Value = "Eval-value: 0.0`nCPSC 313 LEC 01             : MO, 9:00`nCPSC 813 LAB 01             : TU, 18:00"
$ws.Range("B14").WrapText = $true
$ws.Range("C14").ClearContents()
$ws.Rows.Item(14).RowHeight = 45

# --- Row 15: CPSC413 -------------------------------------------------
$ws.Range("B15").Value = "Eval-value: 0.0`nCPSC 413 LEC 01             : MO, 9:00`nCPSC 913 LAB 01             : TU, 18:00"
$ws.Range("B15").WrapText = $true
$ws.Range("C15").ClearContents()
$ws.Rows.Item(15).RowHeight = 45

# --- Row 20: EveningClass --------------------------------------------
$ws.Range("E20").Value = "Eval-value: 0.0`nCPSC 416 LEC 90             : MO, 19:00`nCPSC 518 LEC 91             : MO, 19:00"

# This row is no longer a "problem" case, so drop its yellow highlight.
$ws.Range("A20").Interior.ColorIndex = -4142
$ws.Range("A20").Interior.Pattern = -4142
$ws.Range("B20").Interior.ColorIndex = -4142
$ws.Range("B20").Interior.Pattern = -4142
$ws.Range("E20").Interior.ColorIndex = -4142
$ws.Range("E20").Interior.Pattern = -4142

# The memo explaining the (now-fixed) discrepancy is removed entirely
# (content + formatting), so the cell drops back to the sheet default.
$ws.Range("F20").Clear()

# --- Cosmetic: restore the view/selection as last left by the editor ---
$ws.Activate()
$ws.Range("D20").Select()
